# Updates for test problems
# - species!C2: 0.6 -> 0.2
# - species!C4: 0.1 -> 0.4
# - species!C5: 0.05 -> 0.1
# - species!C6: 0.05 -> 0.1
# - selection on "system" sheet moves to D13
# - selection on "species" sheet moves to C16

$wb = $excel.ActiveWorkbook

$wsSystem = $wb.Worksheets.Item("system")
$wsSpecies = $wb.Worksheets.Item("species")

# Update the recovery/fraction values on the species sheet
$wsSpecies.Range("C2").Value = 0.2
$wsSpecies.Range("C4").Value = 0.4
$wsSpecies.Range("C5").Value = 0.1
$wsSpecies.Range("C6").Value = 0.1

# Update the saved selections for each sheet
$wsSystem.Activate()
$wsSystem.Range("D13").Select()

$wsSpecies.Activate()
$wsSpecies.Range("C16").Select()
